# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# This updates the DAMSLTag (column I) and DialogAct (column J) values for a set of
# rows in the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new DAMSLTag / DialogAct pair
$updates = @{
    8   = @("sv", "Statement-opinion")
    15  = @("sd", "Statement-non-opinion")
    17  = @("sv", "Statement-opinion")
    37  = @("aa", "Agree/Accept")
    38  = @("aa", "Agree/Accept")
    40  = @("aa", "Agree/Accept")
    46  = @("sv", "Statement-opinion")
    56  = @("b",  "Acknowledge (Backchannel)")
    59  = @("sd", "Statement-non-opinion")
    60  = @("ba", "Appreciation")
    61  = @("b",  "Acknowledge (Backchannel)")
    62  = @("sd", "Statement-non-opinion")
    69  = @("sv", "Statement-opinion")
    87  = @("ba", "Appreciation")
    103 = @("b",  "Acknowledge (Backchannel)")
    115 = @("sd", "Statement-non-opinion")
    120 = @("aa", "Agree/Accept")
    134 = @("aa", "Agree/Accept")
    152 = @("b",  "Acknowledge (Backchannel)")
    163 = @("sd", "Statement-non-opinion")
    205 = @("sd", "Statement-non-opinion")
    221 = @("sd", "Statement-non-opinion")
    238 = @("ba", "Appreciation")
    327 = @("aa", "Agree/Accept")
    346 = @("sd", "Statement-non-opinion")
    387 = @("sv", "Statement-opinion")
    389 = @("sv", "Statement-opinion")
    400 = @("sd", "Statement-non-opinion")
    411 = @("sd", "Statement-non-opinion")
    412 = @("sd", "Statement-non-opinion")
    428 = @("sv", "Statement-opinion")
    433 = @("aa", "Agree/Accept")
    435 = @("sd", "Statement-non-opinion")
    441 = @("%",  "Uninterpretable")
    452 = @("sv", "Statement-opinion")
    473 = @("%",  "Uninterpretable")
    484 = @("sd", "Statement-non-opinion")
    493 = @("sv", "Statement-opinion")
    519 = @("aa", "Agree/Accept")
    530 = @("sd", "Statement-non-opinion")
    531 = @("sv", "Statement-opinion")
    533 = @("sv", "Statement-opinion")
    534 = @("sv", "Statement-opinion")
    535 = @("sv", "Statement-opinion")
    545 = @("sd", "Statement-non-opinion")
    547 = @("%",  "Uninterpretable")
    549 = @("b",  "Acknowledge (Backchannel)")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
